$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers (e.g. "225.89") must be
# forced to Text number format first, otherwise Excel auto-converts the entry
# into a numeric value and loses the exact original text (e.g. trailing zeros).
# NumberFormat/Style must be set per-cell (not via a multi-area "A1,A2" Range)
# since multi-area ranges here only apply formatting to the first area.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# --- Apply all cell value updates (document order) ---
$ws.Range("D2").Value = '34.219.52'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.786.49'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '225.89'
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '32.21'
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  -0.54%  '
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").Value = '2.044.66'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").Value = '1.789.46'
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").Value = '10.99'
$ws.Range("E14").Value = '  -4.25%  '
$ws.Range("E15").Value = '  +0.56%  '
$ws.Range("D16").Value = '34.202.70'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '67.97'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("E19").Value = '  +2.42%  '
$ws.Range("D20").Value = '245.93'
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").Value = '10.92'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = '4.16'
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("D25").Value = '161.58'
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("E28").Value = '  +1.67%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").Value = '3.74'
$ws.Range("E32").Value = '  +2.14%  '
$ws.Range("D33").Value = '3.76'
$ws.Range("E33").Value = '  +3.91%  '
$ws.Range("E34").Value = '  -1.56%  '
$ws.Range("D35").Value = '1.436.77'
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("D36").Value = '2.61'
$ws.Range("E36").Value = '  +10.65%  '
$ws.Range("D37").Value = '0.661'
$ws.Range("E37").Value = '  +2.80%  '
$ws.Range("E38").Value = '  +1.30%  '
$ws.Range("E39").Value = '  -1.14%  '
$ws.Range("D40").Value = '81.65'
$ws.Range("E40").Value = '  +1.44%  '
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("D42").Value = '14.08'
$ws.Range("E42").Value = '  +5.58%  '
$ws.Range("E43").Value = '  +1.53%  '
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("E45").Value = '  +1.92%  '
$ws.Range("D46").Value = '6.06'
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  +0.74%  '
$ws.Range("D48").Value = '1.941.09'
$ws.Range("E48").Value = '  -0.44%  '
$ws.Range("D49").Value = '105.30'
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("D51").Value = '0.0₆0130'
$ws.Range("E51").Value = '  -6.38%  '

# Restore the default (unstyled) cell style on the forced-text cells so their
# formatting matches the rest of the sheet, as in the original workbook.
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
